# WorkingHours.xlsx update
# "MockUp, Pflichtenheft und Arbeitszeit auf neustem Stand"
#
# LukasPerger (3rd sheet) gets a new working-hours entry:
#   Datum = 13.03.2024 (serial 45364), Stunden = "2h",
#   Beschreibung = "MockUp verbessert, Projektantrag, Projektauftrag und
#   Pflichtenheft verbessert"
# and becomes the active/selected sheet (instead of PaulSchein).

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)   # LukasPerger

# Note: activating LukasPerger below automatically clears the
# tabSelected flag that was previously on PaulSchein (sheet 1),
# since only one sheet can be the selected tab at a time.

# --- Add the new row of working hours to LukasPerger ---
$ws3.Range("A3").Value = 45364
# Copy the date cell above so the new date cell reuses the same
# date number-format style instead of creating a brand-new one.
$ws3.Range("A2").Copy()
$ws3.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$ws3.Range("B3").Value = "2h"
$ws3.Range("C3").Value = "MockUp verbessert, Projektantrag, Projektauftrag und Pflichtenheft verbessert"

# --- Make LukasPerger the active sheet, with C4 selected ---
$ws3.Activate()
$ws3.Range("C4").Select()
